$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New job posting row (row 12) with Job_Id = JD_011
$row = 12
$ws.Cells.Item($row, 1).Value = "JD_011"
$ws.Cells.Item($row, 2).Value = "Senior SW Engineer"
$ws.Cells.Item($row, 3).Value = "We are seeking a Software Engineer to build and maintain high-quality software solutions.`nWork with global teams to drive innovation and deliver scalable applications.`nJoin Akkodis and be part of a tech-driven, collaborative environment."
$ws.Cells.Item($row, 4).Value = 1
$ws.Cells.Item($row, 5).Value = 4

# Reset row height so it matches the default (no custom height stored)
$ws.Rows.Item($row).AutoFit()

$wb.Save()
